$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new price value is a plain decimal number (e.g. "22.00", "0.9988").
# Excel's COM Value setter auto-converts such text to a number, which would drop
# trailing zeros / switch to scientific notation, so we force Text format first.
$textFormatValues = @{
    "D4" = "0.9988"
    "D6" = "298.47"
    "D7" = "0.3778"
    "D8" = "50.19"
    "D9" = "0.3487"
    "D10" = "0.08050"
    "D11" = "1.209"
    "D12" = "0.9987"
    "D13" = "22.00"
    "D14" = "6.295"
    "D15" = "7.232"
    "D16" = "0.00001196"
    "D18" = "94.71"
    "D19" = "0.06950"
    "D20" = "6.611"
    "D21" = "17.29"
    "D22" = "0.9999"
    "D23" = "12.39"
    "D25" = "2.419"
    "D26" = "2.946"
    "D27" = "20.93"
    "D28" = "149.89"
    "D29" = "5.158"
    "D30" = "131.04"
    "D32" = "6.767"
    "D33" = "2.121"
    "D35" = "0.9867"
    "D36" = "0.02673"
    "D37" = "0.08741"
    "D38" = "0.2415"
    "D39" = "5.854"
    "D40" = "0.06812"
    "D41" = "12.77"
    "D42" = "0.6797"
    "D43" = "1.289"
    "D44" = "15.38"
    "D46" = "0.6322"
    "D47" = "2.228"
    "D48" = "3.891"
    "D49" = "0.07662"
    "D50" = "126.65"
    "D51" = "1.213"
}

foreach ($cellRef in $textFormatValues.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $textFormatValues[$cellRef]
}

# Remaining cells (multi-dot "price" strings and all "Volume(1h)" percentages) are
# never interpreted as numbers by Excel, so they can be assigned directly.
$plainValues = @{
    "D2" = "23.392.31"
    "E2" = "  -1.46%  "
    "D3" = "1.632.67"
    "E3" = "  -1.46%  "
    "E4" = "  -0.21%  "
    "E5" = "  -0.24%  "
    "E6" = "  -1.50%  "
    "E7" = "  -1.41%  "
    "E8" = "  -1.77%  "
    "E9" = "  -3.42%  "
    "E10" = "  -1.78%  "
    "E11" = "  -1.70%  "
    "E12" = "  -0.21%  "
    "E13" = "  -2.05%  "
    "E14" = "  -2.39%  "
    "E15" = "  -2.74%  "
    "E16" = "  -2.33%  "
    "D17" = "1.630.99"
    "E17" = "  -1.54%  "
    "E18" = "  -2.77%  "
    "E19" = "  -1.04%  "
    "E20" = "  -2.75%  "
    "E21" = "  -1.66%  "
    "E22" = "  -0.17%  "
    "E23" = "  -2.74%  "
    "D24" = "23.406.19"
    "E24" = "  -1.42%  "
    "E25" = "  -3.34%  "
    "E26" = "  -2.53%  "
    "E27" = "  -1.46%  "
    "E28" = "  -2.93%  "
    "E29" = "  -1.53%  "
    "E30" = "  -2.23%  "
    "D31" = "1.808.69"
    "E31" = "  -1.71%  "
    "E32" = "  -6.00%  "
    "E33" = "  -5.67%  "
    "E34" = "  -7.17%  "
    "E35" = "  -6.44%  "
    "E36" = "  -4.72%  "
    "E37" = "  -0.65%  "
    "E38" = "  -3.59%  "
    "E39" = "  -4.06%  "
    "E40" = "  -2.58%  "
    "E41" = "  -1.62%  "
    "E42" = "  -2.76%  "
    "E43" = "  -2.90%  "
    "E44" = "  -3.59%  "
    "E45" = "  -0.25%  "
    "E46" = "  -2.92%  "
    "E47" = "  -3.14%  "
    "E48" = "  -1.83%  "
    "E49" = "  -3.01%  "
    "E50" = "  -1.12%  "
    "E51" = "  +1.69%  "
}

foreach ($cellRef in $plainValues.Keys) {
    $ws.Range($cellRef).Value = $plainValues[$cellRef]
}
